# Insert a new record row at row 57, shifting existing rows 57-120 down to 58-121.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("57:57").Insert()

$ws.Range("A57").Value2 = 5
$ws.Range("B57").Value2 = "Macroferia Regional de Talca"
$ws.Range("C57").Value2 = "Maule"
$ws.Range("D57").Value2 = 44894
$ws.Range("E57").Value2 = 7
$ws.Range("F57").Value2 = 100112022
$ws.Range("G57").Value2 = "Arveja Verde"
$ws.Range("H57").Value2 = "Sin especificar"
$ws.Range("I57").Value2 = "Primera"
$ws.Range("J57").Value2 = 400
$ws.Range("K57").Value2 = 19000
$ws.Range("L57").Value2 = 20000
$ws.Range("M57").Value2 = 19500
$ws.Range("N57").Value2 = "`$/saco 25 kilos"
$ws.Range("O57").Value2 = "Región del Maule"
$ws.Range("P57").Value2 = 780
$ws.Range("Q57").Value2 = 25
$ws.Range("R57").Value2 = "Hortaliza"
